$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing row 13 value (45 -> 43)
$ws.Range("A13").Value = 43

# New data to append starting at row 14 (Empenho, Ano Empenho)
$data = @(
    @(45, 2024),
    @(1499, 2024),
    @(3008, 2024),
    @(4005, 2024),
    @(4021, 2024),
    @(6926, 2024),
    @(7108, 2024),
    @(7109, 2024),
    @(8657, 2024),
    @(8684, 2024)
)

$row = 14
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

# Resize the table to include new rows
$table = $ws.ListObjects.Item("EmpenhosTerceirizacao")
$table.Resize($ws.Range("A1:B23"))

# Update selection to match target state (active cell A24)
$ws.Range("A24").Select()
